{"js": "// tiny correction to documentation\n//\n// Net textual changes applied (run-splitting / proofErr / bookmark churn in\n// the original diff is just Word's live-editing noise and is not\n// reproduced):\n//   1. \"Place Config.txt file to the root folder...\" -> an extra space was\n//      left behind after the re-typed word (\"Config.txt<TWO SPACES>file\"),\n//      and \"(here is example of the file)\" gained \" for ODB\" before the\n//      closing parenthesis.\n//   2. \"...recorded CAN messages back to bus.\" -> \"...back to the bus.\"\n//   3. The whole paragraph \"- Playback mode will be switched back\n//      automatically to logger operation after transmission of the last\n//      message from file.\" is removed.\n//   4. \"...LED will toggle every time when message has been sent,\" ->\n//      \"...when the message has been sent,\"\n\nconst body = context.document.body;\n\n// 1a) Re-insert the text right after the bold \"Config.txt\" run so that run's\n//     bold formatting is left untouched; this naturally leaves the extra\n//     space the author's edit introduced.\nlet r1 = body.search(\"file to the root folder\", { matchCase: true });\nr1.load(\"items\");\nawait context.sync();\nif (r1.items.length !== 1) {\n  throw new Error(\"expected 1 match for 'file to the root folder', got \" + r1.items.length);\n}\nr1.items[0].insertText(\" file to the root folder\", \"Replace\");\nawait context.sync();\n\n// 1b) Insert \" for ODB\" right before the closing parenthesis.\nlet r2 = body.search(\"of the file)\", { matchCase: true });\nr2.load(\"items\");\nawait context.sync();\nif (r2.items.length !== 1) {\n  throw new Error(\"expected 1 match for 'of the file)', got \" + r2.items.length);\n}\nr2.items[0].insertText(\"of the file for ODB)\", \"Replace\");\nawait context.sync();\n\n// 2) Insert \"the \" before \"bus.\" in the playback sentence.\nlet r3 = body.search(\"back to bus.\", { matchCase: true });\nr3.load(\"items\");\nawait context.sync();\nif (r3.items.length !== 1) {\n  throw new Error(\"expected 1 match for 'back to bus.', got \" + r3.items.length);\n}\nr3.items[0].insertText(\"back to the bus.\", \"Replace\");\nawait context.sync();\n\n// 3) Delete the standalone \"- Playback mode ...\" paragraph entirely.\nlet r4 = body.search(\n  \"- Playback mode will be switched back automatically to logger operation after transmission of the last message from file.\",\n  { matchCase: true }\n);\nr4.load(\"paragraphs\");\nawait context.sync();\nif (r4.items.length !== 1) {\n  throw new Error(\"expected 1 match for the Playback-mode paragraph, got \" + r4.items.length);\n}\nr4.items[0].paragraphs.getFirst().delete();\nawait context.sync();\n\n// 4) Insert \"the \" before \"message has been sent\".\nlet r5 = body.search(\"when message has been sent\", { matchCase: true });\nr5.load(\"items\");\nawait context.sync();\nif (r5.items.length !== 1) {\n  throw new Error(\"expected 1 match for 'when message has been sent', got \" + r5.items.length);\n}\nr5.items[0].insertText(\"when the message has been sent\", \"Replace\");\nawait context.sync();\n", "ps1": "# tiny correction to documentation\n# Applies the same net textual changes described by the diff:\n#  1. \"...(here is example of the file).\" -> \"...(here is example of the file for ODB).\"\n#  2. \"...recorded CAN messages back to bus.\" -> \"...recorded CAN messages back to the bus.\"\n#  3. Remove the paragraph \"- Playback mode will be switched back automatically to logger\n#     operation after transmission of the last message from file.\" entirely.\n#  4. \"...LED will toggle every time when message has been sent,\" -> \"...LED will toggle\n#     every time when the message has been sent,\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($doc, [string]$searchText, [string]$replaceText) {\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $searchText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.MatchCase = $true\n    $ok = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $ok) {\n        throw \"Replace-FirstMatch: could not find '\" + $searchText + \"'\"\n    }\n}\n\n# 1a) The edit re-typed \"file\" after the bold \"Config.txt\", leaving a double\n#     space behind. Search only the non-bold tail so the \"Config.txt\" run\n#     (and its bold formatting) is left untouched.\nReplace-FirstMatch $d \"file to the root folder\" \" file to the root folder\"\n\n# 1b) Insert \" for ODB\" right before the closing parenthesis.\nReplace-FirstMatch $d \"of the file)\" \"of the file for ODB)\"\n\n# 2) Insert \"the \" before \"bus.\" in the playback sentence.\nReplace-FirstMatch $d \"back to bus.\" \"back to the bus.\"\n\n# 3) Delete the standalone \"- Playback mode ...\" paragraph (text + its paragraph mark).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"- Playback mode will be switched back automatically to logger operation after transmission of the last message from file.\"\n$find.Forward = $true\n$find.Wrap = 0\nif ($find.Execute()) {\n    $matchRange = $find.Parent\n    $para = $matchRange.Paragraphs(1)\n    $para.Range.Delete()\n} else {\n    throw \"could not find the 'Playback mode' paragraph to delete\"\n}\n\n# 4) Insert \"the \" before \"message has been sent\".\nReplace-FirstMatch $d \"when message has been sent\" \"when the message has been sent\"\n"}
